$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Tuesday hours for week of 17 (F17): 3 -> 3.5
$ws.Range("F17").Value = 3.5

# Thursday hours for week of 17 (H17) now logged: 7
$ws.Range("H17").Value = 7

# Move/update the active selection on the sheet to O17 (matches saved view state)
$excel.Goto($ws.Range("O17"))
